# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2  = 8309
    3  = 7733
    9  = 113
    10 = 159
    11 = 227
    12 = 700
    13 = 124
    14 = 1288
    15 = 60
    16 = 49
    17 = 10
    19 = 112
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
